# Ajuste na base enviar
# Replace the phone-number placeholders in column B with a prompt text,
# and move the active selection to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Digite seu número"
$ws.Range("B3").Value = "Digite seu número"

$ws.Range("B4").Select()
